# Adds a new "Tier" column (E) to the schedule sheet, populating a tier
# value for every person row, then leaves the selection where the author
# left it (G17) per the recorded edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("E1").Value = "Tier"

# Per-row tier classification (row 2 .. row 23)
$tiers = @(
    "Tier 1", # row 2  - Arun
    "Tier 1", # row 3  - Christo
    "Tier 2", # row 4  - Amal
    "Tier 2", # row 5  - Rishitha
    "Tier 2", # row 6  - Ridhi
    "Tier 1", # row 7  - Sneha
    "Tier 3", # row 8  - Gurudeep
    "Tier 2", # row 9  - Person11
    "Tier 3", # row 10 - person12
    "Tier 3", # row 11 - Person12
    "Tier 1", # row 12 - person13
    "Tier 1", # row 13 - Person13
    "Tier 2", # row 14 - person14
    "Tier 3", # row 15 - Person14
    "Tier 3", # row 16 - person15
    "Tier 3", # row 17 - Person15
    "Tier 3", # row 18 - person16
    "Tier 2", # row 19 - Person16
    "Tier 2", # row 20 - person17
    "Tier 2", # row 21 - Person17
    "Tier 1", # row 22 - person18
    "Tier 2"  # row 23 - Person18
)

for ($i = 0; $i -lt $tiers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $tiers[$i]
}

# Match the author's final selection/active cell.
$ws.Range("G17").Select()
